# EPBDS-7754 Rename date() to Date()
# Date() is more constructor-like and more friendly for understanding for BAs

$wb = $excel.ActiveWorkbook

# Rename the "date" worksheet tab to "Date"
$dateSheet = $wb.Worksheets.Item("date")
$dateSheet.Name = "Date"

# Update the example code snippet text from "date(...)" to "Date(...)"
$cell = $dateSheet.Range("B5")
$cell.Value = "return Date(year, month, day);"
